# "Grouped response into 'fall','rise' and 'plateau' groups."
#
# Before: one ramp group with columns
#   C=Resting Rate, D=Number of trials, E=Start of ramp, F=End of ramp,
#   G=Start of ramp, H=End of ramp, I=Start of ramp, J=End of ramp
# (three "Start/End of ramp" pairs for up to 3 additional ramps, but no
#  per-group trial count).
#
# After: three full groups, each with its own "Number of trials", "Start of
# ramp" and "End of ramp": C/D/.. becomes
#   C=Number of trials, D=Resting Rate, E=Start of ramp, F=End of ramp,
#   G=Number of trials, H=Start of ramp, I=End of ramp,
#   J=Number of trials, K=Start of ramp, L=End of ramp
# The "Number of trials" value that used to live in D moves to C (grouped
# with the first ramp), D is zeroed out (Resting Rate is no longer tracked
# per-group), and two new "Number of trials" columns (G and J) are
# inserted - one before each of the extra Start/End-of-ramp pairs - filled
# with 0 for every row that has that ramp group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns: one that becomes the "Number of trials"
# column for the 2nd ramp group (before the old column G), and one for the
# 3rd ramp group (before the old column I, which by the time we get there
# has already been shifted right by one to column J).
$ws.Columns.Item(7).Insert() | Out-Null
$ws.Columns.Item(10).Insert() | Out-Null

# --- Header row ---
$ws.Range("C1").Value() = "Number of trials"
$ws.Range("D1").Value() = "Resting Rate"
$ws.Range("G1").Value() = "Number of trials"
$ws.Range("H1").Value() = "Start of ramp"
$ws.Range("I1").Value() = "End of ramp"
$ws.Range("J1").Value() = "Number of trials"
$ws.Range("K1").Value() = "Start of ramp"
$ws.Range("L1").Value() = "End of ramp"

# --- Move the "Number of trials" values from D into C, zero out D ---
for ($r = 2; $r -le 8; $r++) {
    $trials = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value() = $trials
    $ws.Cells.Item($r, 4).Value() = 0
}

# --- Fill the two new "Number of trials" columns with 0 for every row
#     that actually has that ramp group (i.e. already had Start/End of
#     ramp values there before the insert) ---
$ws.Range("G3").Value() = 0
$ws.Range("G4").Value() = 0
$ws.Range("G5").Value() = 0
$ws.Range("G6").Value() = 0
$ws.Range("G7").Value() = 0
$ws.Range("G8").Value() = 0

$ws.Range("J3").Value() = 0
$ws.Range("J4").Value() = 0
$ws.Range("J6").Value() = 0
$ws.Range("J7").Value() = 0
$ws.Range("J8").Value() = 0

# --- Restore the selection to what it was left at in the edited file ---
$ws.Range("L11").Select() | Out-Null
